# Generate Report for handoff
#
# The 92b3dd88-06b7-4bbf-acd7-ac7f81d3b112 file finished its localization
# handback cycle and dropped off the status report; the 57b8156c file moved
# from "Handed back: in sync with en-US" to "Ready for handoff" with a
# refreshed handoff timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Rows("3").Delete()

# --- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-01-11 07:58:05"
$ws.Rows("3").Delete()

# --- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-01-11 07:58:22"
$ws.Rows("3").Delete()
